# "Data source corrected and updated"
# The J/K columns held a mix of a shared-string header ("r"/"s" on row 1)
# and numeric 0.5/1 values for the remaining rows. The corrected data
# source uses a single constant per column for every row (J=1, K=0.5),
# so the header-row strings are no longer needed and the column values
# for rows 2-51 swap too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite J1:J51 and K1:K51 with plain numeric constants (this also
# clears the two leftover shared-string entries "s"/"r" that used to sit
# in J1/K1, since the whole column is now homogeneous numeric data).
$ws.Range("J1:J51").Value = 1
$ws.Range("K1:K51").Value = 0.5

# Reflect the new review position/selection: the author scrolled down to
# row 31 and selected the full K column range that was just corrected.
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("K1:K51").Select() | Out-Null
